$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.017.66'
$ws.Range("E2").Value = '  -2.27%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.102.90'
$ws.Range("E3").Value = '  -1.06%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  -0.81%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '346.88'
$ws.Range("E5").Value = '  +2.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  -0.83%  '
$ws.Range("E7").Value = '  -2.33%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4430'
$ws.Range("E8").Value = '  -3.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09453'
$ws.Range("E9").Value = '  +3.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '52.56'
$ws.Range("E10").Value = '  -3.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.176'
$ws.Range("E11").Value = '  -0.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25.32'
$ws.Range("E12").Value = '  +2.96%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.105.19'
$ws.Range("E13").Value = '  -1.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.752'
$ws.Range("E14").Value = '  -1.78%  '
$ws.Range("E15").Value = '  +0.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '99.81'
$ws.Range("E16").Value = '  +2.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001166'
$ws.Range("E17").Value = '  -1.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.005'
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '20.76'
$ws.Range("E19").Value = '  +5.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06687'
$ws.Range("E20").Value = '  -0.32%  '
$ws.Range("E21").Value = '  -0.90%  '
$ws.Range("E22").Value = '  -4.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.107.59'
$ws.Range("E23").Value = '  -2.18%  '
$ws.Range("E24").Value = '  -2.69%  '
$ws.Range("E25").Value = '  -1.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.349.01'
$ws.Range("E26").Value = '  -1.14%  '
$ws.Range("E27").Value = '  -2.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.556'
$ws.Range("E28").Value = '  -0.29%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '163.12'
$ws.Range("E29").Value = '  -1.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.60'
$ws.Range("E30").Value = '  -1.63%  '
$ws.Range("E31").Value = '  -3.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1060'
$ws.Range("E32").Value = '  -1.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.642'
$ws.Range("E33").Value = '  -1.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.248'
$ws.Range("E34").Value = '  -2.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.959'
$ws.Range("E35").Value = '  +0.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.212'
$ws.Range("E36").Value = '  +4.32%  '
$ws.Range("E37").Value = '  -4.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02575'
$ws.Range("E38").Value = '  -3.89%  '
$ws.Range("E39").Value = '  -1.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2293'
$ws.Range("E40").Value = '  -1.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.58'
$ws.Range("E41").Value = '  -0.81%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6958'
$ws.Range("E42").Value = '  +0.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.316'
$ws.Range("E43").Value = '  +3.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6685'
$ws.Range("E44").Value = '  +2.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.24'
$ws.Range("E45").Value = '  -6.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.303'
$ws.Range("E46").Value = '  -0.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.644'
$ws.Range("E47").Value = '  -1.64%  '
$ws.Range("E48").Value = '  -5.98%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.223'
$ws.Range("E49").Value = '  -3.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '82.37'
$ws.Range("E50").Value = '  -1.97%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07211'
$ws.Range("E51").Value = '  -1.40%  '
